# Refresh cryptos list: updated prices / 1h volume percentages, plus the
# WEMIXToken <-> InternetComputer(DFINITY) rows (34/35) swapping places.
# Leading "'" on Price (column D) values forces text storage so values such
# as "4.44" or "1.00" aren't silently coerced into numbers by Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.213.40"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").Value = "'2.023.02"
$ws.Range("E3").Value = "  -0.08%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'228.43"
$ws.Range("E5").Value = "  +0.63%  "
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'56.03"
$ws.Range("E8").Value = "  +1.59%  "
$ws.Range("D9").Value = "'0.377"
$ws.Range("E9").Value = "  -1.19%  "
$ws.Range("D10").Value = "'0.0781"
$ws.Range("E10").Value = "  -1.62%  "
$ws.Range("D11").Value = "'0.102"
$ws.Range("E11").Value = "  -2.00%  "
$ws.Range("D12").Value = "'2.322.83"
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("D13").Value = "'14.28"
$ws.Range("E13").Value = "  -0.14%  "
$ws.Range("D14").Value = "'20.19"
$ws.Range("E14").Value = "  -2.22%  "
$ws.Range("D15").Value = "'0.740"
$ws.Range("E15").Value = "  -0.92%  "
$ws.Range("E16").Value = "  +0.58%  "
$ws.Range("D17").Value = "'2.020.58"
$ws.Range("E17").Value = "  -0.03%  "
$ws.Range("D18").Value = "'37.167.56"
$ws.Range("E18").Value = "  +0.48%  "
$ws.Range("D19").Value = "'6.15"
$ws.Range("E19").Value = "  +1.81%  "
$ws.Range("D20").Value = "'68.88"
$ws.Range("E20").Value = "  -0.08%  "
$ws.Range("D21").Value = "'0.0₃0817"
$ws.Range("E21").Value = "  -1.71%  "
$ws.Range("D22").Value = "'223.12"
$ws.Range("E22").Value = "  -1.24%  "
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("E24").Value = "  +1.81%  "
$ws.Range("E25").Value = "  -1.63%  "
$ws.Range("D26").Value = "'163.53"
$ws.Range("E26").Value = "  -2.39%  "
$ws.Range("D27").Value = "'9.02"
$ws.Range("E27").Value = "  -3.71%  "
$ws.Range("D28").Value = "'0.130"
$ws.Range("E28").Value = "  +2.73%  "
$ws.Range("D29").Value = "'18.69"
$ws.Range("E29").Value = "  -0.57%  "
$ws.Range("E31").Value = "  +0.26%  "
$ws.Range("E32").Value = "  -1.02%  "
$ws.Range("D33").Value = "'0.0604"
$ws.Range("E33").Value = "  -1.26%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").Value = "'4.44"
$ws.Range("E34").Value = "  -0.46%  "
$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").Value = "'1.99"
$ws.Range("E35").Value = "  +9.08%  "
$ws.Range("D36").Value = "'2.32"
$ws.Range("E36").Value = "  -2.25%  "
$ws.Range("E37").Value = "  +0.49%  "
$ws.Range("D38").Value = "'1.00"
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("D39").Value = "'5.56"
$ws.Range("E39").Value = "  +2.07%  "
$ws.Range("D40").Value = "'4.40"
$ws.Range("E40").Value = "  +19.06%  "
$ws.Range("D41").Value = "'1.466.97"
$ws.Range("E41").Value = "  -2.41%  "
$ws.Range("E42").Value = "  -3.10%  "
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("D44").Value = "'93.87"
$ws.Range("E44").Value = "  -1.56%  "
$ws.Range("E45").Value = "  -1.55%  "
$ws.Range("D46").Value = "'16.21"
$ws.Range("E46").Value = "  -5.46%  "
$ws.Range("E47").Value = "  -2.23%  "
$ws.Range("E48").Value = "  +0.20%  "
$ws.Range("D49").Value = "'7.13"
$ws.Range("E49").Value = "  -1.63%  "
$ws.Range("D50").Value = "'2.92"
$ws.Range("E50").Value = "  +0.91%  "
$ws.Range("D51").Value = "'2.212.11"
$ws.Range("E51").Value = "  +0.04%  "
